# FINFLUX-2544  Automating Nabkisan Sanity Scenario
#
# Re-pastes the recalculated loan-schedule figures (interest/principal
# recomputed after a sanity-scenario tweak) into the verification
# workbook, and nudges the recorded cell selections on a couple of
# sheets to match where the author was last clicked.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Summary" sheet — totals reflecting the recalculated schedule
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("F2").Value = 3269.11
$summary.Range("A3").Value = 352.83
$summary.Range("E3").Value = 352.83
$summary.Range("F3").Value = 183.72

# ---------------------------------------------------------------------
# "Repayment Schedule" sheet — per-installment breakdown + header tweak
# ---------------------------------------------------------------------
$sched = $wb.Worksheets.Item("Repayment Schedule")

# The "In Advance" header slides from column M to column N (column M
# is left blank); "Late" stays put in column O.
$sched.Range("N1").Value = "In Advance"
$sched.Range("M1").ClearContents()

# Row 3
$sched.Range("F3").NumberFormat = "#,##0.00"
$sched.Range("F3").Value = 1627.83
$sched.Range("G3").NumberFormat = "#,##0.00"
$sched.Range("G3").Value = 8372.17
$sched.Range("K3").NumberFormat = "#,##0.00"
$sched.Range("K3").Value = 1727.83
$sched.Range("N3").Value = 0
$sched.Range("M3").ClearContents()
$sched.Range("Q3").NumberFormat = "#,##0.00"
$sched.Range("Q3").Value = 1727.83

# Row 4
$sched.Range("F4").Value = 1641.28
$sched.Range("G4").Value = 6730.89
$sched.Range("H4").Value = 83.72
$sched.Range("N4").Value = 0
$sched.Range("M4").ClearContents()

# Row 5
$sched.Range("F5").Value = 1657.69
$sched.Range("G5").Value = 5073.2
$sched.Range("H5").Value = 67.31
$sched.Range("N5").Value = 0
$sched.Range("M5").ClearContents()

# Row 6
$sched.Range("F6").Value = 1674.27
$sched.Range("G6").Value = 3398.93
$sched.Range("H6").Value = 50.73
$sched.Range("N6").Value = 0
$sched.Range("M6").ClearContents()

# Row 7
$sched.Range("F7").Value = 1691.01
$sched.Range("G7").Value = 1707.92
$sched.Range("H7").Value = 33.99
$sched.Range("N7").Value = 0
$sched.Range("M7").ClearContents()

# Row 8
$sched.Range("F8").Value = 1707.92
$sched.Range("H8").Value = 17.08
$sched.Range("N8").Value = 0
$sched.Range("M8").ClearContents()

# ---------------------------------------------------------------------
# Recorded selections: author had last clicked J13 on the Repayment
# Schedule sheet. Re-select Summary's own remembered cell afterwards so
# the active tab/selection reported for Summary is undisturbed.
# ---------------------------------------------------------------------
$sched.Range("J13").Select()
$summary.Activate()
$summary.Range("C8").Select()
